$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.603.53"
$ws.Range("E2").Value = "  +0.44%  "
$ws.Range("D3").Value = "2.117.97"
$ws.Range("E3").Value = "  +1.26%  "
$ws.Range("E4").Value = "  +0.71%  "
$ws.Range("D5").Value = "336.74"
$ws.Range("E5").Value = "  +1.99%  "
$ws.Range("E6").Value = "  +0.64%  "
$ws.Range("D7").Value = "0.5250"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").Value = "0.4552"
$ws.Range("E8").Value = "  +3.39%  "
$ws.Range("D9").Value = "54.59"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "0.09134"
$ws.Range("E10").Value = "  +2.37%  "
$ws.Range("D11").Value = "1.174"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").Value = "24.49"
$ws.Range("E12").Value = "  +1.19%  "
$ws.Range("D13").Value = "2.119.63"
$ws.Range("E13").Value = "  +1.07%  "
$ws.Range("D14").Value = "6.859"
$ws.Range("E14").Value = "  +2.72%  "
$ws.Range("D15").Value = "8.154"
$ws.Range("E15").Value = "  +6.09%  "
$ws.Range("D16").Value = "0.00001177"
$ws.Range("E16").Value = "  +4.96%  "
$ws.Range("D17").Value = "97.11"
$ws.Range("E17").Value = "  +1.31%  "
$ws.Range("D18").Value = "1.010"
$ws.Range("E18").Value = "  +0.71%  "
$ws.Range("D19").Value = "0.06675"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  +1.58%  "
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "6.309"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("D23").Value = "30.664.53"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("D24").Value = "12.92"
$ws.Range("E24").Value = "  +5.13%  "
$ws.Range("D25").Value = "2.352"
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "2.363.87"
$ws.Range("E26").Value = "  +0.86%  "
$ws.Range("D27").Value = "22.43"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "2.557"
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "134.52"
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("E31").Value = "  +1.96%  "
$ws.Range("E32").Value = "  +0.44%  "
$ws.Range("E33").Value = "  -0.06%  "
$ws.Range("D34").Value = "6.371"
$ws.Range("E34").Value = "  +3.66%  "
$ws.Range("D35").Value = "3.944"
$ws.Range("E35").Value = "  +1.17%  "
$ws.Range("D36").Value = "10.66"
$ws.Range("E36").Value = "  +5.74%  "
$ws.Range("D37").Value = "5.879"
$ws.Range("E37").Value = "  +7.56%  "
$ws.Range("D38").Value = "0.02630"
$ws.Range("E38").Value = "  +2.98%  "
$ws.Range("D39").Value = "0.06845"
$ws.Range("E39").Value = "  +0.61%  "
$ws.Range("D41").Value = "12.60"
$ws.Range("E41").Value = "  +0.11%  "
$ws.Range("D42").Value = "0.6896"
$ws.Range("E42").Value = "  +0.51%  "
$ws.Range("D43").Value = "1.258"
$ws.Range("E43").Value = "  +0.91%  "
$ws.Range("D44").Value = "14.85"
$ws.Range("E44").Value = "  +6.53%  "
$ws.Range("D45").Value = "0.6494"
$ws.Range("E45").Value = "  +2.87%  "
$ws.Range("D46").Value = "2.311"
$ws.Range("E46").Value = "  +5.47%  "
$ws.Range("D47").Value = "0.00000000365"
$ws.Range("E47").Value = "  +20.90%  "
$ws.Range("E48").Value = "  +1.72%  "
$ws.Range("E49").Value = "  +1.01%  "
$ws.Range("D50").Value = "83.44"
$ws.Range("E50").Value = "  +2.35%  "
$ws.Range("B51").Value = "WEMIXTOKEN"
$ws.Range("C51").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D51").Value = "1.189"
$ws.Range("E51").Value = "  -3.88%  "
